$d = $word.ActiveDocument
$endash = [char]0x2013

function Get-ParagraphIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.StartsWith($text)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1. "Tony - 8 hours" -> "Tony - 40 hours"
# ---------------------------------------------------------------------
$idxTony = Get-ParagraphIndexByText $d "Tony"
$rTony = $d.Paragraphs($idxTony).Range
$rTony.Find.Execute("8 hours", $true, $false, $false, $false, $false, $true, 1, $false, "40 hours", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Insert two new bullet paragraphs after
#    "Add code functionality to each user - 4 hours - Risk: 2"
# ---------------------------------------------------------------------
$idxAddCode = Get-ParagraphIndexByText $d "Add code functionality"
$d.Paragraphs($idxAddCode).Range.InsertParagraphAfter()

$idxFix = $idxAddCode + 1
$d.Paragraphs($idxFix).Range.Text = "Fix some documentation " + $endash + " 2 hours " + $endash + " Risk: 1"

$d.Paragraphs($idxFix).Range.InsertParagraphAfter()
$idxAvail = $idxFix + 1
$d.Paragraphs($idxAvail).Range.Text = "Available time to help others " + $endash + " 34 hours " + $endash + " Risk: 2"

# ---------------------------------------------------------------------
# 3. "...managing employees - 10 hours  - Risk: 3" ->
#    "...managing employees - 10 hours - Risk: 3"   (drop the doubled space)
# ---------------------------------------------------------------------
$idxManaging = Get-ParagraphIndexByText $d "Include usage of database in managing"
if ($idxManaging -eq -1) {
    $idxManaging = Get-ParagraphIndexByText $d "Include usage of database in m"
}
$rManaging = $d.Paragraphs($idxManaging).Range
$rManaging.Find.Execute("hours  " + $endash, $true, $false, $false, $false, $false, $true, 1, $false, "hours " + $endash, 2) | Out-Null

Write-Host "done"
